$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The original edit landed on shape id 6 ("TextBox 5") even though only
# ids 2-4 were in use beforehand; adding + deleting a throwaway shape
# first reproduces that same id/name allocation.
$tmp = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$tmp.Delete()

$shp = $s.Shapes.AddTextbox(1, 255.7636220472441, 9.963779527559055, 480.10897637795273, 29.081259842519685)

$shp.TextFrame.WordWrap = -1
$shp.TextFrame.AutoSize = 1
$shp.Fill.Visible = 0

$tr = $shp.TextFrame.TextRange
$tr.Text = "https://github.com/OneOfTheInfiniteMonkeys/MTMP"
$tr.LanguageID = "en-GB"
